$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 data
$ws.Range("A2").Value = "281474991109864-1743437700258"
$ws.Range("B2").Value = "Mobile Usage"
$ws.Range("C2").Value = "2025-03-31T10:15:00.258"
$ws.Range("D2").Value = "'281474991109864"
$ws.Range("E2").Value = "'138"
$ws.Range("F2").Value = "'51833957"
$ws.Range("G2").Value = "FERNANDO ORNELAS"
$ws.Range("H2").Value = 20.66838345
$ws.Range("I2").Value = -103.418986399
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743437697758/qrwRQlnS6o-camera-video-segment-driver-1743437700258.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIASGIR6T%2F20250401%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250401T151230Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEFUaCXVzLXdlc3QtMiJGMEQCIDWsZf5%2BbVu7V2d55jueBNIwQA3BkyBQfdhsjqN4ExyoAiBz6JNn04%2Bn1qQ%2FEwxMhQIyjzQ1%2BINcxOQSyJgJl4eN2irmAwi9%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMqVhkufHVPlk9ACPAKroDcYcgduFEcegZiRCdzlou6cVWG6ZaN0gds0YnZphPbCMPuP%2BebXHIM5xTyjeBeUD%2B5%2BTLZUr2sRYYKnijx5jvcwSkkv12TooDHaZ98iF7GWlf7XQhL9jlCsJF2zJa6lMbGqquvK03Lo0LvSdlTMkZkCUsCGG9S1acr%2B%2BrAcl1%2F682pLVloVWp%2B4xab3igpnvdIK%2BZsTRz1oOxCjfuyhk6GSCYziDsd4G%2FNF%2F1Na1dls8vVaNzWSAJYSZ7nleERf%2FC43hSBPBwksye3UqSI3ZAU86TdYnJ5WoiR2WFeYLOU%2B%2FtYByzRqunWgId40kSCXeGPWRBY5mafVEGbRl9PUrxID%2Fs%2BtLRLc2lJZTUfqDfp21NK0HXb8Mpb910AdAsPFnmstBEbqnQw9Kd4oQH3QJI4INlMTd9X%2FrhC5yuQE54WdpuO1TJs%2FBxcUckv09faLo7NvZ0QuXSRtyo9BKim7DOZbZSQvHIuzT5%2BnbtRHQcYJY8HzP45ce%2BBkD9k1DHh0CRSDHhQLovg8%2BMGmAtaOKdg9hhDqBgru1kTyiJAC1x5o0OTO6%2BJLzmezQDJZOzr%2F0el0h%2FYqAGm8RwpzDmta%2B%2FBjqmAcm5mgPnZGbkP4hgZqaZI8pmWtQsxc4NEaM%2FxtUqD0vCRJGN585uokrCTW9KKAQGJDEaZzCAVHc%2BR5GsIqzaBr3sxuKk3C8GS6t7XWr2FEmhYmnr1rGrZsM%2FYc8oQsxOuChLLn7qxLeuIL15M1CNJCA1GOlNxkUTxJ4ZVxmUlmsL2%2FpfxBsHaNhW2%2BKPPMj4zia2PHmGqjBMmMyEeOv9e2pr7ldcYY0%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2001%20Apr%202025%2023%3A12%3A30%20GMT&X-Amz-Signature=b259f432a80b969c4af11862793c9fd6a8e52ab5991ff01105c4a5a92d5bd0ca'
$ws.Range("L2").Value = "No video URL"

# Row 3 data
$ws.Range("A3").Value = "281474991109864-1743433748258"
$ws.Range("B3").Value = "Mobile Usage"
$ws.Range("C3").Value = "2025-03-31T09:09:08.258"
$ws.Range("D3").Value = "'281474991109864"
$ws.Range("E3").Value = "'138"
$ws.Range("F3").Value = "'51833957"
$ws.Range("G3").Value = "FERNANDO ORNELAS"
$ws.Range("H3").Value = 20.66134081
$ws.Range("I3").Value = -103.386754229
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743433745758/xb4UbhC61W-camera-video-segment-driver-1743433748258.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIASGIR6T%2F20250401%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250401T151230Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEFUaCXVzLXdlc3QtMiJGMEQCIDWsZf5%2BbVu7V2d55jueBNIwQA3BkyBQfdhsjqN4ExyoAiBz6JNn04%2Bn1qQ%2FEwxMhQIyjzQ1%2BINcxOQSyJgJl4eN2irmAwi9%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMqVhkufHVPlk9ACPAKroDcYcgduFEcegZiRCdzlou6cVWG6ZaN0gds0YnZphPbCMPuP%2BebXHIM5xTyjeBeUD%2B5%2BTLZUr2sRYYKnijx5jvcwSkkv12TooDHaZ98iF7GWlf7XQhL9jlCsJF2zJa6lMbGqquvK03Lo0LvSdlTMkZkCUsCGG9S1acr%2B%2BrAcl1%2F682pLVloVWp%2B4xab3igpnvdIK%2BZsTRz1oOxCjfuyhk6GSCYziDsd4G%2FNF%2F1Na1dls8vVaNzWSAJYSZ7nleERf%2FC43hSBPBwksye3UqSI3ZAU86TdYnJ5WoiR2WFeYLOU%2B%2FtYByzRqunWgId40kSCXeGPWRBY5mafVEGbRl9PUrxID%2Fs%2BtLRLc2lJZTUfqDfp21NK0HXb8Mpb910AdAsPFnmstBEbqnQw9Kd4oQH3QJI4INlMTd9X%2FrhC5yuQE54WdpuO1TJs%2FBxcUckv09faLo7NvZ0QuXSRtyo9BKim7DOZbZSQvHIuzT5%2BnbtRHQcYJY8HzP45ce%2BBkD9k1DHh0CRSDHhQLovg8%2BMGmAtaOKdg9hhDqBgru1kTyiJAC1x5o0OTO6%2BJLzmezQDJZOzr%2F0el0h%2FYqAGm8RwpzDmta%2B%2FBjqmAcm5mgPnZGbkP4hgZqaZI8pmWtQsxc4NEaM%2FxtUqD0vCRJGN585uokrCTW9KKAQGJDEaZzCAVHc%2BR5GsIqzaBr3sxuKk3C8GS6t7XWr2FEmhYmnr1rGrZsM%2FYc8oQsxOuChLLn7qxLeuIL15M1CNJCA1GOlNxkUTxJ4ZVxmUlmsL2%2FpfxBsHaNhW2%2BKPPMj4zia2PHmGqjBMmMyEeOv9e2pr7ldcYY0%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2001%20Apr%202025%2023%3A12%3A30%20GMT&X-Amz-Signature=0cff63853adea50d3c8a10fbb87745fd7f0759a6ccf576d022941761fcf13eaf'
$ws.Range("L3").Value = "No video URL"
